$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 7.371418320197964
$ws.Range("G2").Value = 2864
$ws.Range("H2").Value = 5.461644439755614
$ws.Range("I2").Value = 1.349670122525919
$ws.Range("K2").Value = 198.4015080113101
$ws.Range("L2").Value = 1083.598493069101
$ws.Range("M2").Value = 802.8617326440752
$ws.Range("F3").Value = 7.42844776884932
$ws.Range("G3").Value = 2870
$ws.Range("H3").Value = 5.492392392159672
$ws.Range("I3").Value = 1.352497643732328
$ws.Range("K3").Value = 170.4147031102733
$ws.Range("L3").Value = 935.9844188750144
$ws.Range("M3").Value = 692.0414414121186
$ws.Range("F4").Value = 7.317252841704606
$ws.Range("G4").Value = 2870
$ws.Range("H4").Value = 5.41017788505128
$ws.Range("I4").Value = 1.352497643732328
$ws.Range("K4").Value = 1.352497643732328
$ws.Range("L4").Value = 7.317252841704606
$ws.Range("M4").Value = 5.41017788505128
$ws.Range("F5").Value = 7.363302375792443
$ws.Range("G5").Value = 2868
$ws.Range("H5").Value = 5.44802219017837
$ws.Range("I5").Value = 1.351555136663525
$ws.Range("K5").Value = 1.351555136663525
$ws.Range("L5").Value = 7.363302375792443
$ws.Range("M5").Value = 5.44802219017837
$ws.Range("F6").Value = 8.31199416601871
$ws.Range("G6").Value = 2876
$ws.Range("H6").Value = 6.132841314426879
$ws.Range("I6").Value = 1.355325164938737
$ws.Range("K6").Value = 178.9029217719133
$ws.Range("L6").Value = 1097.18322991447
$ws.Range("M6").Value = 809.5350535043481
$ws.Range("F7").Value = 8.501961772048368
$ws.Range("G7").Value = 2923
$ws.Range("H7").Value = 6.17213919955068
$ws.Range("I7").Value = 1.377474081055608
$ws.Range("K7").Value = 202.4886899151743
$ws.Range("L7").Value = 1249.78838049111
$ws.Range("M7").Value = 907.3044623339499
$ws.Range("F8").Value = 8.414732436597731
$ws.Range("G8").Value = 2924
$ws.Range("H8").Value = 6.10672442902202
$ws.Range("I8").Value = 1.377945334590009
$ws.Range("K8").Value = 1.377945334590009
$ws.Range("L8").Value = 8.414732436597731
$ws.Range("M8").Value = 6.10672442902202
$ws.Range("F9").Value = 8.339491127254979
$ws.Range("G9").Value = 2874
$ws.Range("H9").Value = 6.157411333345534
$ws.Range("I9").Value = 1.354382657869934
$ws.Range("K9").Value = 106.9962299717248
$ws.Range("L9").Value = 658.8197990531434
$ws.Range("M9").Value = 486.4354953342972
$ws.Range("F10").Value = 10.63626539872964
$ws.Range("G10").Value = 2932
$ws.Range("H10").Value = 7.697870114633113
$ws.Range("I10").Value = 1.381715362865221
$ws.Range("K10").Value = 1.381715362865221
$ws.Range("L10").Value = 10.63626539872964
$ws.Range("M10").Value = 7.697870114633113
$ws.Range("F11").Value = 10.48732659047415
$ws.Range("G11").Value = 2875
$ws.Range("H11").Value = 7.740558965212569
$ws.Range("I11").Value = 1.354853911404336
$ws.Range("K11").Value = 102.9688972667295
$ws.Range("L11").Value = 797.0368208760351
$ws.Range("M11").Value = 588.2824813561552
$ws.Range("F12").Value = 10.84743884553431
$ws.Range("G12").Value = 3196
$ws.Range("H12").Value = 7.202210647754633
$ws.Range("I12").Value = 1.50612629594722
$ws.Range("K12").Value = 118.9839773798303
$ws.Range("L12").Value = 856.9476687972106
$ws.Range("M12").Value = 568.974641172616
$ws.Range("F13").Value = 10.8950341608467
$ws.Range("G13").Value = 3205
$ws.Range("H13").Value = 7.213498436604276
$ws.Range("I13").Value = 1.510367577756833
$ws.Range("K13").Value = 114.7879359095193
$ws.Range("L13").Value = 828.0225962243495
$ws.Range("M13").Value = 548.2258811819249
$ws.Range("F14").Value = 11.74624383926375
$ws.Range("G14").Value = 3444
$ws.Range("H14").Value = 7.237377882380278
$ws.Range("I14").Value = 1.622997172478793
$ws.Range("K14").Value = 1.622997172478793
$ws.Range("L14").Value = 11.74624383926375
$ws.Range("M14").Value = 7.237377882380278
$ws.Range("F15").Value = 10.24147735122762
$ws.Range("G15").Value = 3478
$ws.Range("H15").Value = 6.248537935395345
$ws.Range("I15").Value = 1.639019792648445
$ws.Range("K15").Value = 129.4825636192271
$ws.Range("L15").Value = 809.076710746982
$ws.Range("M15").Value = 493.6344968962322
$ws.Range("F16").Value = 10.61233894270542
$ws.Range("G16").Value = 3574
$ws.Range("H16").Value = 6.300890664919108
$ws.Range("I16").Value = 1.68426013195099
$ws.Range("K16").Value = 1.68426013195099
$ws.Range("L16").Value = 10.61233894270542
$ws.Range("M16").Value = 6.300890664919108
$ws.Range("F17").Value = 10.28120131816876
$ws.Range("G17").Value = 3483
$ws.Range("H17").Value = 6.263769508226848
$ws.Range("I17").Value = 1.641376060320452
$ws.Range("K17").Value = 124.7445805843544
$ws.Range("L17").Value = 781.3713001808259
$ws.Range("M17").Value = 476.0464826252405
$ws.Range("F18").Value = 10.384868231269
$ws.Range("G18").Value = 3488
$ws.Range("H18").Value = 6.317858482440599
$ws.Range("I18").Value = 1.64373232799246
$ws.Range("K18").Value = 1.64373232799246
$ws.Range("L18").Value = 10.384868231269
$ws.Range("M18").Value = 6.317858482440599
$ws.Range("F19").Value = 11.4159614464458
$ws.Range("G19").Value = 3781
$ws.Range("H19").Value = 6.4069479474631
$ws.Range("I19").Value = 1.781809613572102
$ws.Range("K19").Value = 1.781809613572102
$ws.Range("L19").Value = 11.4159614464458
$ws.Range("M19").Value = 6.4069479474631
$ws.Range("F20").Value = 12.39498951216258
$ws.Range("G20").Value = 3575
$ws.Range("H20").Value = 7.357249718827689
$ws.Range("I20").Value = 1.684731385485391
$ws.Range("K20").Value = 133.0937794533459
$ws.Range("L20").Value = 979.2041714608437
$ws.Range("M20").Value = 581.2227277873875
$ws.Range("F21").Value = 12.62355898840928
$ws.Range("G21").Value = 3640
$ws.Range("H21").Value = 7.359118728957278
$ws.Range("I21").Value = 1.715362865221489
$ws.Range("K21").Value = 1.715362865221489
$ws.Range("L21").Value = 12.62355898840928
$ws.Range("M21").Value = 7.359118728957278
$ws.Range("F22").Value = 12.41449192126643
$ws.Range("G22").Value = 3571
$ws.Range("H22").Value = 7.377079769512006
$ws.Range("I22").Value = 1.682846371347785
$ws.Range("K22").Value = 1.682846371347785
$ws.Range("L22").Value = 12.41449192126643
$ws.Range("M22").Value = 7.377079769512006
$ws.Range("F23").Value = 12.40397839045271
$ws.Range("G23").Value = 3573
$ws.Range("H23").Value = 7.366706449633544
$ws.Range("I23").Value = 1.683788878416588
$ws.Range("K23").Value = 127.9679547596607
$ws.Range("L23").Value = 942.7023576744061
$ws.Range("M23").Value = 559.8696901721494
$ws.Range("F24").Value = 13.55117374673068
$ws.Range("G24").Value = 3841
$ws.Range("H24").Value = 7.486485470076155
$ws.Range("I24").Value = 1.810084825636192
$ws.Range("K24").Value = 1.810084825636192
$ws.Range("L24").Value = 13.55117374673068
$ws.Range("M24").Value = 7.486485470076155
$ws.Range("K25").Value = 1.459542306457408
$ws.Range("L25").Value = 9.413538577795778
$ws.Range("M25").Value = 6.449651055778068
$ws.Range("N25").Value = 40.22637508256916
$ws.Range("O25").Value = 2122
$ws.Range("Q25").Value = 0.02163130269970093

$ws.Range("R25").Value = "(38.52088324656055, 41.931866918577775)"
$ws.Range("S25").Value = "(37.98487152667213, 42.467878638466196)"

Write-Output "done"
